# Re-create the excel sheet's data cells: both A1 and B1 collapse down to a
# single blank (space) value instead of the "Data"/"Service" header labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = " "
$ws.Range("B1").Value = " "
